# Update "想去人数" (column F) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row => new value for column F
$updates = @{
    2  = 11827
    3  = 11567
    6  = 1044
    8  = 73
    9  = 47
    11 = 10882
    12 = 4199
    19 = 6
    20 = 136
    22 = 11163
    23 = 10970
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
